# Apply the "resultados" sheet restructuring for arboltd, bosquetd and knntd:
#   - remove the "Arreglo aleatorio óptimo <sheet>" / seed column (old column B)
#     so MAE/MSE/RMSE/R2 slide one column to the left (B..E instead of C..F)
#   - drop the per-sheet suffix from the first header ("... arboltd" -> "...")
#   - rename "R-cuadrado <sheet>" -> "R2 <sheet>"
#   - refresh the row-2 metric values with the newly computed results

$wb = $excel.ActiveWorkbook

function Update-ResultSheet {
    param($SheetName, $HeaderA, $Suffix, $A2, $B2, $C2, $D2, $E2)

    $ws = $wb.Worksheets.Item($SheetName)

    # Overwrite B..E in place with the (shifted) MAE/MSE/RMSE/R2 columns,
    # then fully clear the now-unused F column (contents + formatting) so
    # it drops out of the sheet and the dimension shrinks to E2. We
    # deliberately avoid a structural column delete/insert here: this sheet
    # is referenced by formulas on other sheets (e.g. ResltNumericas), and a
    # real column delete would rewrite those formulas (...!B2, ...!F2, etc.)
    # — which the target edit does not do.
    $ws.Range("A1").Value = $HeaderA
    $ws.Range("B1").Value = "MAE $Suffix"
    $ws.Range("C1").Value = "MSE $Suffix"
    $ws.Range("D1").Value = "RMSE $Suffix"
    $ws.Range("E1").Value = "R2 $Suffix"

    $ws.Range("A2").Value = $A2
    $ws.Range("B2").Value = $B2
    $ws.Range("C2").Value = $C2
    $ws.Range("D2").Value = $D2
    $ws.Range("E2").Value = $E2

    $ws.Range("F1:F2").Clear()
}

Update-ResultSheet "bosquetd" "Estimador óptimo" "bosquetd" 122 0.4573004169495208 0.410033964531995 0.6403389450377004 0.6226367865212324

Update-ResultSheet "knntd" "K óptimo" "knntd" 15 0.6220102972122369 0.6818153889718599 0.8257211326906075 0.3725104053871119

Update-ResultSheet "arboltd" "Profundidad óptima" "arboltd" 4 0.6845854597204523 0.7617854841998375 0.8728032333807189 0.2989121800501523
